# Weekly data refresh: insert a new daily/weekly record as row 394,
# pushing all the subsequent existing records down by one row
# (old row 394 becomes row 395, ..., old row 442 becomes row 443).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 394; Excel shifts rows 394:442 down to 395:443
# and extends the used range to row 443.
$ws.Rows("394:394").Insert()

# Populate the newly inserted row 394 with the new weekly observation.
$ws.Range("A394").Value2 = 10
$ws.Range("B394").Value2 = "Vega Modelo de Temuco"
$ws.Range("C394").Value2 = "La Araucanía"
$ws.Range("D394").Value2 = 44946
$ws.Range("E394").Value2 = 9
$ws.Range("F394").Value2 = 100112009
$ws.Range("G394").Value2 = "Acelga"
$ws.Range("H394").Value2 = "Sin especificar"
$ws.Range("I394").Value2 = "Primera"
$ws.Range("J394").Value2 = 35
$ws.Range("K394").Value2 = 8000
$ws.Range("L394").Value2 = 8000
$ws.Range("M394").Value2 = 8000
$ws.Range("N394").Value2 = "`$/docena de atados (12 kilos)"
$ws.Range("O394").Value2 = "Provincia de Cautín"
$ws.Range("P394").Value2 = 667
$ws.Range("Q394").Value2 = 12
$ws.Range("R394").Value2 = "Hortaliza"
